$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.067.80"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.651.65"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D5").Value = "'217.24"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'0.5262"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D8").Value = "'0.2598"
$ws.Range("D9").Value = "'0.06326"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'20.36"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "'0.07796"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'4.504"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "1.651.19"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "0.0₅8206"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "'65.53"
$ws.Range("D17").Value = "26.071.88"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'4.577"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "'190.92"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "'10.07"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'6.023"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'143.46"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'7.223"
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "'15.99"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "'1.428"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "'0.05816"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "'1.272"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'3.551"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "'3.273"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'1.580"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'0.9462"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.781"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").Value = "'0.5735"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'0.01610"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.8429"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'5.743"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'103.65"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "1.028.98"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("D44").Value = "1.795.65"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "'56.83"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'0.4322"
$ws.Range("E47").Value = "  +2.92%  "
$ws.Range("D48").Value = "'7.872"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'1.463"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'0.09621"
$ws.Range("E51").Value = "  -0.90%  "
